$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label rows 2-5 to reflect the reordered metric names
# (PackageCycle, Clique, Crossing, UnhealthyInheritance)
$ws.Range("A2").Value = "PackageCycle"
$ws.Range("A3").Value = "Clique"
$ws.Range("A4").Value = "Crossing"
$ws.Range("A5").Value = "UnhealthyInheritance"

# Update the numeric data (columns B..F) for each row with the new measurements
$data = @{
    2 = @(437.0, 328.0, 6089.0, 4392.0, 190886.0)
    3 = @(557.0, 399.0, 7451.0, 5994.0, 258196.0)
    4 = @(626.0, 489.0, 9239.0, 7409.0, 314094.0)
    5 = @(661.0, 501.0, 10814.0, 6943.0, 299387.0)
    6 = @(730.0, 511.0, 9361.0, 8126.0, 337219.0)
    7 = @(933.0, 548.0, 10080.0, 9145.0, 379858.0)
    9 = @(1215.0, 600.0, 13541.0, 9765.0, 406702.0)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
}
